$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Remember the width of the column immediately to the left (M) so the
# newly inserted column can inherit the same formatting/width, which is
# what Excel does automatically when a column is inserted from the UI.
$leftWidth = $ws.Columns("M").ColumnWidth

# Insert a new blank column before column N ("Late"), shifting the
# existing N/O/P ("Late", "Date"/heading, "Outstanding") columns one to
# the right (-> O/P/Q).
$ws.Columns("N").Insert()

# The new column N inherits the left neighbour's width (as Excel itself
# would do), instead of the default column width.
$ws.Columns("N").ColumnWidth = $leftWidth

# Make "Repayment schedule" the active sheet/tab and leave the selection
# where the user ended up after the edit.
$ws.Activate()
$ws.Range("L15").Select()
